# Natmi following Dr Hou advice
# Lpl-Vldlr LR-pair results recomputed: every Sending-cluster x Target-cluster
# combination (ECs / FAPs / sCs) is now reported -> 9 data rows instead of 6.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

$data = @(
  ,@("ECs", "Lpl", "Vldlr", "ECs", 3, 1, 23.496322, 70.488966, 0.1321353991144917, 0.1321353991144917, 1, 0.3333333333333333, 0.06089466666666667, 0.182684, 0.001903591634475228, 0.001903591634475228, 1.430800696082667, 12.877206264744, 0.0002515318403723918, 0.0002515318403723917)
  ,@("ECs", "Lpl", "Vldlr", "FAPs", 3, 1, 23.496322, 70.488966, 0.1321353991144917, 0.1321353991144917, 3, 1, 28.046323, 84.138969, 0.8767392739472014, 0.8767392739472013, 658.9854361240061, 5930.868925116055, 0.1158482938823631, 0.1158482938823631)
  ,@("ECs", "Lpl", "Vldlr", "sCs", 3, 1, 23.496322, 70.488966, 0.1321353991144917, 0.1321353991144917, 3, 1, 3.882136333333333, 11.646409, 0.1213571344183235, 0.1213571344183235, 91.21592533589934, 820.9433280230941, 0.01603557339175618, 0.01603557339175618)
  ,@("FAPs", "Lpl", "Vldlr", "ECs", 3, 1, 124.9120333333333, 374.7361, 0.7024631931202969, 0.7024631931202969, 1, 0.3333333333333333, 0.06089466666666667, 0.182684, 0.001903591634475228, 0.001903591634475228, 7.606476632488889, 68.4582896924, 0.001337203057950554, 0.001337203057950554)
  ,@("FAPs", "Lpl", "Vldlr", "FAPs", 3, 1, 124.9120333333333, 374.7361, 0.7024631931202969, 0.7024631931202969, 3, 1, 28.046323, 84.138969, 0.8767392739472014, 0.8767392739472013, 3503.323233453433, 31529.9091010809, 0.6158770699109218, 0.6158770699109217)
  ,@("FAPs", "Lpl", "Vldlr", "sCs", 3, 1, 124.9120333333333, 374.7361, 0.7024631931202969, 0.7024631931202969, 3, 1, 3.882136333333333, 11.646409, 0.1213571344183235, 0.1213571344183235, 484.9255430738777, 4364.3298876649, 0.08524892015142459, 0.08524892015142459)
  ,@("sCs", "Lpl", "Vldlr", "ECs", 3, 1, 29.411685, 88.235055, 0.1654014077652114, 0.1654014077652114, 1, 0.3333333333333333, 0.06089466666666667, 0.182684, 0.001903591634475228, 0.001903591634475228, 1.79101475418, 16.11913278762, 0.0003148567361522824, 0.0003148567361522824)
  ,@("sCs", "Lpl", "Vldlr", "FAPs", 3, 1, 29.411685, 88.235055, 0.1654014077652114, 0.1654014077652114, 3, 1, 28.046323, 84.138969, 0.8767392739472014, 0.8767392739472013, 824.889617484255, 7424.006557358295, 0.1450139101539164, 0.1450139101539164)
  ,@("sCs", "Lpl", "Vldlr", "sCs", 3, 1, 29.411685, 88.235055, 0.1654014077652114, 0.1654014077652114, 3, 1, 3.882136333333333, 11.646409, 0.1213571344183235, 0.1213571344183235, 114.180170963055, 1027.621538667495, 0.02007264087514269, 0.02007264087514269)
)

$startRow = 2
for ($i = 0; $i -lt $data.Length; $i++) {
  $r = $startRow + $i
  $rowVals = $data[$i]
  for ($j = 0; $j -lt $cols.Length; $j++) {
    $ws.Cells.Item($r, $j + 1).Value = $rowVals[$j]
  }
}
